$d = $word.ActiveDocument

# The "first page" header (header3.xml) contains the signature block with
# the doctor's name. Update it to add the "PD" prefix.
$wdHeaderFooterFirstPage = 2

$section = $d.Sections.Item(1)
$hdr = $section.Headers.Item($wdHeaderFooterFirstPage)

$hdr.Range.Find.Execute("Dr. med. Thiên-Trí Lâm", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PD Dr. med. Thiên-Trí Lâm", 2)
